# davie trainingsplan 2014 - "made initial data transformation"
#
# Restructures the weekly training-plan sheet:
#  - adds a "Datum" header label in column A
#  - regroups the per-exercise "set N" columns so each exercise's 4 sets
#    are contiguous (Planke, Liegestuetz, Hammer Curls, Bizeps Curls,
#    Kniebeugen, Turm Rudern, Shoulder Press)
#  - renames the weight-lifting exercises to "Weighted ..." (Hammer Curls,
#    Bizeps Curls, Turm Rudern, Shoulder Press), with the multi-word ones
#    wrapped onto their own header lines
#  - header row gets wrapped text + taller row height to fit the new labels
#  - the weight|reps notation in the data rows switches from "2, NN" to a
#    "2 | NN" pipe-separated style
#  - columns B:AC get a narrower, uniform width

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1: header labels, column by column, in the new layout
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datum"

$ws.Range("B1").Value = "Planke set 1"
$ws.Range("C1").Value = "Planke set 2"
$ws.Range("D1").Value = "Planke set 3"
$ws.Range("E1").Value = "Planke set 4"

$ws.Range("F1").Value = "Liegestütz set 1"
$ws.Range("G1").Value = "Liegestütz set 2"
$ws.Range("H1").Value = "Liegestütz set 3"
$ws.Range("I1").Value = "Liegestütz set 4"

$ws.Range("J1").Value = "Weighted Hammer Curls set 1"
$ws.Range("K1").Value = "Weighted Hammer Curls set 2"
$ws.Range("L1").Value = "Weighted Hammer Curls set 3"
$ws.Range("M1").Value = "Weighted Hammer Curls set 4"

$ws.Range("N1").Value = "Weighted `nBizeps Curls `nset 1"
$ws.Range("O1").Value = "Weighted `nBizeps Curls `nset 2"
$ws.Range("P1").Value = "Weighted `nBizeps Curls `nset 3"
$ws.Range("Q1").Value = "Weighted `nBizeps Curls `nset 4"

$ws.Range("R1").Value = "Kniebeugen set 1"
$ws.Range("S1").Value = "Kniebeugen set 2"
$ws.Range("T1").Value = "Kniebeugen set 3"
$ws.Range("U1").Value = "Kniebeugen set 4"

$ws.Range("V1").Value = "Weighted `nTurm Rudern `nset 1"
$ws.Range("W1").Value = "Weighted `nTurm Rudern `nset 2"
$ws.Range("X1").Value = "Weighted `nTurm Rudern `nset 3"
$ws.Range("Y1").Value = "Weighted Turm Rudern `nset 4"

$ws.Range("Z1").Value = "Weighted Shoulder Press set 1"
$ws.Range("AA1").Value = "Weighted Shoulder Press set 2"
$ws.Range("AB1").Value = "Weighted Shoulder Press set 3"
$ws.Range("AC1").Value = "Weighted Shoulder Press set 4"

# ---------------------------------------------------------------------
# Row 1 formatting: wrap the new multi-line labels and grow the row
# ---------------------------------------------------------------------
$headerRow = $ws.Range("A1:AC1")
$headerRow.WrapText = $true
$ws.Rows.Item(1).RowHeight = 45

# ---------------------------------------------------------------------
# Data rows: weight|reps notation "2, NN" -> "2 | NN" (pipe separated)
# ---------------------------------------------------------------------
$ws.Range("J2").Value = "2| 12"
$ws.Range("K2").Value = "2 |15"

$ws.Range("J5").Value = "2 | 18"
$ws.Range("K5").Value = "2 | 14"

$ws.Range("J7").Value = "2 | 10"
$ws.Range("K7").Value = "2 | 18"

# ---------------------------------------------------------------------
# Columns B:AC get a narrower, uniform width (was 16.77.../18.77...)
# ---------------------------------------------------------------------
$ws.Range("B1:AC1").ColumnWidth = 13.95

# ---------------------------------------------------------------------
# Restore the cursor/selection to where the editor left off
# ---------------------------------------------------------------------
$ws.Range("M8").Select()
